$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Oxyspheraster euaster"
$ws.Range("A2").Value = "forked oxyaster"
[void]$ws.Range("A3").Select()
